# "Same problem solved in R" - re-solved the Airline Revenue Management
# linear program with a different solver/engine, producing slightly
# different optimal values on Sheet1 (the model sheet). The cached
# Solver report sheets (Answer/Sensitivity/Limits Report) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Inputs updated by the new solve:
#  - Regular Demand (D5) dropped from 200 to 150
#  - Regular Seats decision (E5) now 150 (was 166)
#  - Discount Seats decision (E6) now 100 (was 0)
#  - Capacity constraint RHS (D11) raised from 166 to 250
$ws.Range("D5").Value = 150
$ws.Range("E5").Value = 150
$ws.Range("E6").Value = 100
$ws.Range("D11").Value = 250

$wb.Application.Calculate()
